# Guarding Pawns translation sheet update.
#
# The translation source moved on to a new export ("Main_250701") while the
# previous export ("Main_240430") is kept around for reference, renamed to
# "Old_240430" and pushed to the second tab. A new row documenting the
# melee-attack job report string is only added to the new sheet.

$wb = $excel.ActiveWorkbook

# A throwaway sheet bumps the internal sheetId counter so the duplicated
# sheet below lands on sheetId 3 (matching the authored workbook), then it
# is discarded.
$scratch = $wb.Worksheets.Add()

$wb.Worksheets.Item("Main_240430").Copy($wb.Worksheets.Item("Main_240430"), $null)

$scratch.Delete()

$wb.Worksheets.Item("Main_240430 (2)").Name = "Main_250701"
$wb.Worksheets.Item("Main_240430").Name = "Old_240430"

$newSheet = $wb.Worksheets.Item("Main_250701")
$oldSheet = $wb.Worksheets.Item("Old_240430")

# New row describing the melee-attack job report string.
$newSheet.Range("A93").Value = "JobDef+GuardingP_AttackMelee.reportString"
$newSheet.Range("B93").Value = "JobDef"
$newSheet.Range("C93").Value = "GuardingP_AttackMelee.reportString"
$newSheet.Range("F93").Value = "TargetA 공격 중"

# Selection/view bookkeeping to mirror the authored workbook state.
$oldSheet.Cells.Select()
$newSheet.Select()
$newSheet.Range("A80").Select()
